# assign6.docx: turn the Mondak & Sanders citation into a hyperlink.
#
# Before (one bullet, plain text):
#   "Replication: Read Mondak and Sanders (2003) - "Tolerance and
#   Intolerance, 1976-1998"."
#
# After:
#   "Replication: Browse " + [hyperlink: "Mondak and Sanders (2003) -
#   "Tolerance and Intolerance, 1976-1998""] + "."

$d = $word.ActiveDocument

# 1) "Read Mondak" -> "Browse Mondak" (shrinks the lead-in text so it now
#    reads "Replication: Browse" immediately before the linked text).
$r1 = $d.Content
$found1 = $r1.Find.Execute("Read Mondak", $true, $false, $false, $false, $false, $true, 1, $false, "Browse Mondak", 2)
if (-not $found1) {
    throw "Could not find 'Read Mondak' to rewrite as 'Browse Mondak'."
}

# 2) Locate "Mondak and Sanders (2003) - "Tolerance and Intolerance,
#    1976-1998"" and wrap it in a hyperlink to the article.
$linkText = "Mondak and Sanders (2003) - " + [char]0x201C + "Tolerance and Intolerance, 1976-1998" + [char]0x201D

$r2 = $d.Content
$found2 = $r2.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the Mondak and Sanders citation to hyperlink."
}

$null = $d.Hyperlinks.Add($r2, "https://www.jstor.org/stable/3186116", "", "", "")
Write-Host "Hyperlinked the Mondak and Sanders citation."
